$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For D-column (price) values that look numeric, force Text number format before
# assigning the string, then clear formats afterward so the cell keeps the default
# (unstyled) appearance while the stored value stays a literal text string (e.g. "1.00").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.683.01"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.36%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.324.48"
$ws.Range("D3").ClearFormats()

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.99"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.49"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.539"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  +8.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.29"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.74%  "

$ws.Range("E11").Value = "  +3.77%  "

$ws.Range("E12").Value = "  -0.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.07"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.683.73"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.12"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.322.24"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.818"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.576.61"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.60"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0939"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.59%  "

$ws.Range("E21").Value = "  +2.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.55"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.41"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.06"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.64"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  -1.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.86"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.89"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.62%  "

$ws.Range("E30").Value = "  +3.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.72"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.92"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.70%  "

$ws.Range("E33").Value = "  +2.29%  "

$ws.Range("E34").Value = "  +0.02%  "

# Row 35/36: coin entries swap places (LidoDAOToken now ranks above WEMIXToken)
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.14"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.25%  "

$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.53"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.78%  "

$ws.Range("E37").Value = "  +1.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.84"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.88%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.108"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.59%  "

$ws.Range("E40").Value = "  +2.97%  "

$ws.Range("E41").Value = "  +2.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.36"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +7.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.01"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +6.13%  "

$ws.Range("E44").Value = "  -0.27%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0292"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.986.58"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.09%  "

$ws.Range("E47").Value = "  +5.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.93"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "56.09"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.62%  "

$ws.Range("E50").Value = "  +2.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.58"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +7.58%  "
